$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "v_cmd_idx := shared_cmd_idx",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "v_cmd_idx := get_last_received_cmd_idx(SBI_VVCT, 1)", 2
)
